$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Length of reach (m)" column between the existing "GNIS Name"
# column (A) and the "Distance (m)" column (old B, which shifts to C, while
# old C "Apportionment (%)" shifts to D).
# Shift the header row (6) and the data row (7) one column to the right,
# working right-to-left so each copy's source isn't clobbered before it's read.
$ws.Range("C6").Copy($ws.Range("D6"))
$ws.Range("B6").Copy($ws.Range("C6"))
$ws.Range("C7").Copy($ws.Range("D7"))
$ws.Range("B7").Copy($ws.Range("C7"))

# Fill in the new column's header + template placeholder cells.
$ws.Range("B6").Value = "Length of reach (m)"
$ws.Range("B7").Value = "{d.streams[i].length_metre}"

# The title-row fill used to span B1:D1; now that column B carries real
# header text, drop the extra formatted-but-empty C1/D1 cells.
$ws.Range("C1").Clear()
$ws.Range("D1").Clear()

# Column widths: widen B for the new header text, and size the new column D.
$ws.Columns("B").ColumnWidth = 18.17
$ws.Columns("D").ColumnWidth = 21.33

# Move the active selection to B8 (matches the saved view state).
[void]$ws.Range("B8").Select()

Write-Host "ok"
